$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 4 (WoodenClub) previously used a pair of redundant "applyFont" style
# variants; restyle it to match the plain styles used by the rest of the
# table (same visible formatting as row 3).
$ws.Range("A3:E3").Copy()
$ws.Range("A4:E4").PasteSpecial(-4122)

# Fill in the new item row: ID=9, NAME=UpgradeI, TYPE=USABLE, TIER=NONE, ATTRIBUTES=NONE
$ws.Range("A11").Value = 9
$ws.Range("B11").Value = "UpgradeI"
$ws.Range("C11").Value = "USABLE"
$ws.Range("D11").Value = "NONE"
$ws.Range("E11").Value = "NONE"

# Copy formatting from the row above (row 10) so the new row matches the table style
$ws.Range("A10:E10").Copy()
$ws.Range("A11:E11").PasteSpecial(-4122)

# Move the active selection to E12, where the next entry would go
$ws.Range("E12").Select()
